# Add a "Save" column (H) to the s_vals sheet, mirroring the style of the
# existing header row and filling in the per-row Save flag values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell H1 - copy style from existing header (e.g. G1) so it matches
# the bold/bordered/centered header formatting, then set its text.
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Per-row "Save" values for rows 2-21 (matches H2:H21 in the target sheet).
$saveValues = @(0,0,0,0,0,1,0,0,0,0,0,0,1,0,0,0,0,0,0,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
